# Practical 8 (autophagosome) - fill in the "count" column (D) with the
# actual per-GO-term counts computed during the tutorial. Most rows keep
# their original "0" count; a handful of rows get updated counts.
# The values must be written as *text* (matching the existing "0" entries,
# which are also stored as text/shared-strings rather than numbers), so we
# enter them with a leading apostrophe and then restore the cell's default
# "Normal" style so no extra number formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$counts = @{
    "D3"  = "1"
    "D4"  = "10"
    "D5"  = "3"
    "D7"  = "2"
    "D15" = "1"
    "D24" = "1"
    "D28" = "2"
    "D33" = "5"
    "D36" = "2"
}

foreach ($addr in $counts.Keys) {
    $rng = $ws.Range($addr)
    $rng.Value = "'" + $counts[$addr]
    $rng.Style = "Normal"
}
